$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The sheet currently has 10 data rows (2-10) plus a totals row (11).
# The edit: update rows 2-10 content, insert two brand-new data rows
# (the old totals row shifts from 11 down to 13), and refresh the
# totals row's numeric values to match the new data.
# ---------------------------------------------------------------------

# Insert two blank rows before the current totals row (row 11) so it
# shifts down to row 13, making room for new rows 11 and 12.
$ws.Rows.Item(11).Resize(2).Insert()

# r, A, B, C, D, E, F, G, H, I, J, K, L, M, N, O
$data = @(
    @(2,  "389/AOURIR/AV1", "Point de vente", "FF", "AGENCE KHATABI", "oui", "mensuelle", 0, 3000, 0, 0, 0, 0, 0, "--", 3000),
    @(3,  "389/AOURIR/AV1", "Point de vente", "A6743213", "ZEROUALI IBTISSAM", "non", "mensuelle", 10, 4500, 0, 450, 0, 0, 0, "--", 4050),
    @(4,  "001/SUP SUD", "Supervision", "1098777", "AGENCE LAHLOU", "oui", "mensuelle", 0, 3000, 0, 0, 0, 0, 0, "--", 3000),
    @(5,  "001/SUP SUD", "Supervision", $null, "BENNIS MOHAMED", "non", "mensuelle", 10, 3000, 0, 300, 0, 0, 0, "--", 2700),
    @(6,  "001/SUP SUD", "Supervision", "B12346", "BAKKALI MOHAMED", "non", "mensuelle", 10, 3000, 0, 300, 0, 0, 0, "--", 2700),
    @(7,  "988/DIRECTION CAPITAL SOFT", "Direction régionale", "B12346", "BAKKALI MOHAMED", "non", "mensuelle", 0, 2000, 0, 0, 0, 0, 0, "--", 2000),
    @(8,  "988/DIRECTION CAPITAL SOFT", "Direction régionale", "A123456", "YOUSSEF", "non", "mensuelle", 10, 4000, 0, 400, 0, 0, 0, "--", 3600),
    @(9,  "988/DIRECTION CAPITAL SOFT", "Direction régionale", "J207703", "ACHENGLI LAILA", "non", "mensuelle", 0, 2000, 0, 0, 0, 0, 0, "--", 2000),
    @(10, "604/ERRAHMA", "Point de vente", "19087", "AGENCE ESSALAM", "oui", "mensuelle", 0, 3333.33, 0, 0, 0, 0, 0, "--", 3333.33),
    @(11, "604/ERRAHMA", "Point de vente", "A6743213", "ZEROUALI IBTISSAM", "non", "mensuelle", 10, 3333.33, 0, 333.33, 0, 0, 0, "--", 3000),
    @(12, "604/ERRAHMA", "Point de vente", $null, "EL OUAZZANI SIHAM", "non", "mensuelle", 10, 3333.33, 0, 333.33, 0, 0, 0, "--", 3000)
)

# Purely-numeric "CIN/IF" codes (column C) must stay text, not become
# numbers - force the text format before writing them.
$ws.Cells.Item(4, 3).NumberFormat = "@"
$ws.Cells.Item(10, 3).NumberFormat = "@"

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    if ($null -ne $row[3]) {
        $ws.Cells.Item($r, 3).Value = $row[3]
    } else {
        $ws.Cells.Item($r, 3).ClearContents()
    }
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
    $ws.Cells.Item($r, 14).Value = $row[14]
    $ws.Cells.Item($r, 15).Value = $row[15]
}

# Update the totals row, now at row 13, to match the new data sums.
$ws.Cells.Item(13, 8).Value = 34499.99
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 2116.66
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(13, 14).Value = 0
$ws.Cells.Item(13, 15).Value = 32383.33
